$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: new note in column H ---
$ws.Range("H4").Value = "Find missing columns for each business_id. Eg business_id : 3,[column1, column2, column3]"

# --- Row 3: task text updated (LA -> states/Las Vegas wording) ---
$ws.Range("B3").Value = "Isolate data according to states and Las Vegas from yelp data - coding"

# --- Row 16: assignee ---
$ws.Range("E16").Value = "Rajdeep"

# --- Row 7: assignee filled in ---
$ws.Range("E7").Value = "Done"

# --- Row 12: assignee ---
$ws.Range("E12").Value = "Rajdeep and Rimsha"

# --- Row 17: assignee ---
$ws.Range("E17").Value = "Devershi!!!???"

# --- Rows 10, 11, 13, 14, 15: assignee / person column filled in ---
$ws.Range("E10").Value = "Ajay"
$ws.Range("E11").Value = "Ajay"
$ws.Range("E13").Value = "Devershi"
$ws.Range("E14").Value = "Devershi"
$ws.Range("E15").Value = "Rimsha"

# --- View/selection housekeeping to mirror the author's saved state ---
$ws.Activate() | Out-Null
$ws.Range("E16").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
